# pilot-04 update: refresh muscle response tStart/tStop values and selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated response-time measurements (tStart = col B, tStop = col C)
$ws.Range("B3").Value  = 0.002
$ws.Range("C3").Value  = 0.0035

$ws.Range("B4").Value  = 0.002
$ws.Range("C4").Value  = 0.0035

$ws.Range("B5").Value  = 0.0022
$ws.Range("C5").Value  = 0.0032

$ws.Range("B6").Value  = 0.0022
$ws.Range("C6").Value  = 0.0032

$ws.Range("B7").Value  = 0.002
$ws.Range("C7").Value  = 0.0045

$ws.Range("B8").Value  = 0.002
$ws.Range("C8").Value  = 0.0045

$ws.Range("B9").Value  = 0.002
$ws.Range("C9").Value  = 0.0045

$ws.Range("B10").Value = 0.002
$ws.Range("C10").Value = 0.0045

$ws.Range("B11").Value = 0.003
$ws.Range("C11").Value = 0.006

$ws.Range("B12").Value = 0.003
$ws.Range("C12").Value = 0.006

$ws.Range("B13").Value = 0.003
$ws.Range("C13").Value = 0.006

$ws.Range("B14").Value = 0.003
$ws.Range("C14").Value = 0.006

$ws.Range("B15").Value = 0.00085
$ws.Range("C15").Value = 0.0045

$ws.Range("B16").Value = 0.00085
$ws.Range("C16").Value = 0.0045

# Move the active selection to B13 (was F13) and scroll back to top of sheet.
$ws.Range("B13").Select()
